# Added scenarios for PUT program
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# Row 2: programDescription gets the new "tester" scenario value
$ws.Range("D2").Value = "tester"

# Row 25: new RubyAPI/TestS PUT-program scenario
$ws.Range("B25").Value = "RubyAPI"
$ws.Range("D25").Value = "TestS"
$ws.Range("E25").Value = "200"

# Row 2: programName becomes the new "Wells" scenario value
$ws.Range("B2").Value = "Wells"

# Row 3: new "PUT program" valid scenario (was the old JavaX/400/Invalid-desc row)
$ws.Range("B3").Value = "Fargo"

# D3's number moves down to D4 (row 3 becomes blank there); restore D4's General
# number format temporarily so the value round-trips as a real number, then put
# the original Text format back so the cell style stays s="11".
$ws.Range("D3").Value = $null
$ws.Range("A1").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = 1234
$ws.Range("D5").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("E3").Value = "201"
$ws.Range("G3").Value = "Valid With Mandatory Details"
$ws.Range("H3").Value = "Program2"

# Row 4: now holds what used to be row 3's scenario (minus the H column note)
$ws.Range("B4").Value = "JavaX"
$ws.Range("E4").Value = "400"
$ws.Range("G4").Value = "Invalid Program Desc"
$ws.Range("H4").Value = $null

# Update the active selection to match the authored workbook
$ws.Range("B3").Select() | Out-Null
